$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the test result values into column C
$ws.Range("C1").Value = "PASSED"
$ws.Range("C2").Value = "PASSED"

# Resize column B to fit its contents (e.g. "Domino's Pizza")
$ws.Columns("B:B").ColumnWidth = 13.42

# Leave the active selection on the last written cell
$ws.Range("C2").Select() | Out-Null
